{"js": "// Replace each \"a\u00f7b=\" division prompt in the practice-sheet table with its\n// updated value. Every old value in the document is unique, so a simple\n// search-and-replace keyed on the exact old text is unambiguous.\nconst replacements = [\n  [\"78\u00f74=\", \"18\u00f76=\"],\n  [\"52\u00f74=\", \"56\u00f79=\"],\n  [\"14\u00f78=\", \"12\u00f77=\"],\n  [\"83\u00f76=\", \"13\u00f78=\"],\n  [\"71\u00f76=\", \"21\u00f73=\"],\n  [\"19\u00f76=\", \"81\u00f73=\"],\n  [\"39\u00f72=\", \"27\u00f72=\"],\n  [\"55\u00f78=\", \"23\u00f73=\"],\n  [\"81\u00f79=\", \"58\u00f73=\"],\n  [\"22\u00f79=\", \"26\u00f75=\"],\n  [\"45\u00f77=\", \"26\u00f74=\"],\n  [\"69\u00f79=\", \"43\u00f77=\"],\n  [\"89\u00f78=\", \"12\u00f73=\"],\n  [\"50\u00f76=\", \"65\u00f76=\"],\n  [\"21\u00f78=\", \"34\u00f74=\"],\n  [\"84\u00f75=\", \"36\u00f77=\"],\n  [\"28\u00f78=\", \"90\u00f74=\"],\n  [\"47\u00f73=\", \"44\u00f75=\"],\n  [\"20\u00f78=\", \"57\u00f74=\"],\n  [\"27\u00f76=\", \"33\u00f78=\"],\n  [\"85\u00f76=\", \"79\u00f75=\"],\n  [\"15\u00f78=\", \"75\u00f72=\"],\n  [\"63\u00f73=\", \"20\u00f79=\"],\n  [\"38\u00f77=\", \"40\u00f79=\"],\n  [\"61\u00f78=\", \"42\u00f72=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace each \"a\u00f7b=\" division prompt in the practice-sheet table with its\n# updated value. Every old value in the document is unique, so a simple\n# Find/Replace keyed on the exact old text is unambiguous.\n$d = $word.ActiveDocument\n\n$pairs = @(\n  @(\"78\u00f74=\", \"18\u00f76=\"),\n  @(\"52\u00f74=\", \"56\u00f79=\"),\n  @(\"14\u00f78=\", \"12\u00f77=\"),\n  @(\"83\u00f76=\", \"13\u00f78=\"),\n  @(\"71\u00f76=\", \"21\u00f73=\"),\n  @(\"19\u00f76=\", \"81\u00f73=\"),\n  @(\"39\u00f72=\", \"27\u00f72=\"),\n  @(\"55\u00f78=\", \"23\u00f73=\"),\n  @(\"81\u00f79=\", \"58\u00f73=\"),\n  @(\"22\u00f79=\", \"26\u00f75=\"),\n  @(\"45\u00f77=\", \"26\u00f74=\"),\n  @(\"69\u00f79=\", \"43\u00f77=\"),\n  @(\"89\u00f78=\", \"12\u00f73=\"),\n  @(\"50\u00f76=\", \"65\u00f76=\"),\n  @(\"21\u00f78=\", \"34\u00f74=\"),\n  @(\"84\u00f75=\", \"36\u00f77=\"),\n  @(\"28\u00f78=\", \"90\u00f74=\"),\n  @(\"47\u00f73=\", \"44\u00f75=\"),\n  @(\"20\u00f78=\", \"57\u00f74=\"),\n  @(\"27\u00f76=\", \"33\u00f78=\"),\n  @(\"85\u00f76=\", \"79\u00f75=\"),\n  @(\"15\u00f78=\", \"75\u00f72=\"),\n  @(\"63\u00f73=\", \"20\u00f79=\"),\n  @(\"38\u00f77=\", \"40\u00f79=\"),\n  @(\"61\u00f78=\", \"42\u00f72=\")\n)\n\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $range = $d.Content\n  $range.Find.Text = $oldText\n  $range.Find.Replacement.Text = $newText\n  $range.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n"}
